# IndividualWorkSheet.xlsx — "Add files via upload" update
#
# The uploaded workbook advances 유병주 (3rd sheet)'s work-log with a new
# week's entries: the old rows 5-8 (Java/MySql integration notes) are
# replaced/shifted so that row 5 keeps the Java-Eclipse/MySql testing entry
# (now paired with a new short task name), and two brand-new log rows are
# added describing (1) reflecting menu-ordered ingredient stock, and
# (2) reworking the Menu file's ingredient input field. The former row 8
# becomes a new blank template row (matching rows 9-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

function Set-TextCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.VerticalAlignment = -4160   # xlTop
    $r.WrapText = $true
    $r.Value = $text
}

function Set-DateCell($addr, $serial) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "yyyy\-mm\-dd;@"
    $r.VerticalAlignment = -4160   # xlTop
    $r.Value = $serial
}

# --- Row 5: keep the Java Eclipse/MySql testing entry, but retitle task (A)
#     and refresh the result/problem notes (E/F). Dates (C/D) stay the same.
Set-TextCell "A5" "Java와 Mysql 연동하기"
Set-TextCell "B5" "Java Eclipse와 MySql DB 연동 및 testing"
Set-DateCell "C5" 43600
Set-DateCell "D5" 43603
Set-TextCell "E5" "testing 코드 git commit"
Set-TextCell "F5" "실제 gui java 파일에서 tomcat server가 연결문제 처리하는데 시간 필요"

# --- Row 6: new log entry - reflecting ingredient stock on menu orders
Set-TextCell "A6" "메뉴 주문 시, 재료 재고량을 반영하기"
Set-TextCell "B6" "테이블에서 메뉴를 주문할 시, 해당 메뉴의 재료가 줄어든다. 재료가 부족하다면 주문할 수 없다."
Set-DateCell "C6" 43594
Set-DateCell "D6" 43597
Set-TextCell "E6" "Table 파일 commit"
Set-TextCell "F6" "보안성이 낮음. 개선이 필요함"

# --- Row 7: new log entry - reworking Menu file's ingredient input field
Set-TextCell "A7" "Menu 파일의 재료 입력란 수정"
Set-TextCell "B7" "기존에는 한 메뉴에 하나의 재료만 입력할 수 있었음. 수정 뒤엔 여러 메뉴 입력 가능"
Set-DateCell "C7" 43594
Set-DateCell "D7" 43597
Set-TextCell "E7" "Menu 파일 commit"
Set-TextCell "F7" "재료 문자열을 분리하는 기능을 넣지 못함. 수정 필요."

# --- Row 8: no longer holds data - clear back to a blank template row
$row8 = $ws.Range("A8:F8")
$row8.ClearContents()
$ws.Rows.Item(8).AutoFit()

# --- Column C:D get a touch wider (best-fit) now that row 6/7 carry dates
$ws.Range("C1:D1").EntireColumn.AutoFit()

# --- Selection/scroll follows the newly-edited rows
$ws.Activate()
$ws.Range("H6").Select()
